$wb = $excel.ActiveWorkbook

# --- "survey" sheet: date/time fields for the JGI app were switched to plain text ---
$survey = $wb.Worksheets.Item("survey")

# Row 2 = OS_FOL_date ("Date of follow"): type column was "date" -> "text"
$survey.Range("C2").Value = "text"

# Row 4 = OS_time_begin ("Begin time of encounter"): type column was "time" -> "text"
$survey.Range("C4").Value = "text"

# Row 5 = OS_time_end ("End time of encounter"): type column was "time" -> "text"
$survey.Range("C5").Value = "text"

# Make "survey" the active sheet/tab and move the selection to C7
$survey.Activate()
$survey.Range("C7").Select()

# --- restore the workbook window to its default (un-moved, un-resized) layout ---
$win = $excel.ActiveWindow
$win.WindowState = -4143  # xlNormal
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 16060
